$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2122186495176849
$ws.Range("C2").Value = 0.5401929260450161
$ws.Range("J2").Value = 0.009646302250803859
$ws.Range("P2").Value = 0.135048231511254
$ws.Range("S2").Value = 0.1028938906752412
$ws.Range("C3").Value = 0.005882352941176471
$ws.Range("J3").Value = 0.06470588235294118
$ws.Range("P3").Value = 0.6823529411764706
$ws.Range("S3").Value = 0.2470588235294118
$ws.Range("J4").Value = 0.1136363636363636
$ws.Range("P4").Value = 0.5454545454545454
$ws.Range("S4").Value = 0.3409090909090909
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.0410958904109589
$ws.Range("D6").Value = 0.0136986301369863
$ws.Range("F6").Value = 0.0502283105022831
$ws.Range("J6").Value = 0.273972602739726
$ws.Range("O6").Value = 0.0273972602739726
$ws.Range("Q6").Value = 0.1141552511415525
$ws.Range("R6").Value = 0.0867579908675799
$ws.Range("S6").Value = 0.3926940639269406
$ws.Range("B7").Value = 0.1337209302325581
$ws.Range("D7").Value = 0.005813953488372093
$ws.Range("F7").Value = 0.05813953488372093
$ws.Range("J7").Value = 0.186046511627907
$ws.Range("O7").Value = 0.02325581395348837
$ws.Range("Q7").Value = 0.1686046511627907
$ws.Range("R7").Value = 0.06395348837209303
$ws.Range("S7").Value = 0.3604651162790697
$ws.Range("B8").Value = 0.1103286384976526
$ws.Range("D8").Value = 0.01643192488262911
$ws.Range("F8").Value = 0.06338028169014084
$ws.Range("J8").Value = 0.136150234741784
$ws.Range("O8").Value = 0.02347417840375587
$ws.Range("Q8").Value = 0.1877934272300469
$ws.Range("R8").Value = 0.08215962441314555
$ws.Range("S8").Value = 0.3802816901408451
$ws.Range("B9").Value = 0.1391304347826087
$ws.Range("D9").Value = 0.01739130434782609
$ws.Range("F9").Value = 0.1043478260869565
$ws.Range("J9").Value = 0.0782608695652174
$ws.Range("O9").Value = 0.01739130434782609
$ws.Range("Q9").Value = 0.1652173913043478
$ws.Range("R9").Value = 0.1043478260869565
$ws.Range("S9").Value = 0.3739130434782609
$ws.Range("B10").Value = 0.1239168110918544
$ws.Range("D10").Value = 0.0268630849220104
$ws.Range("E10").Value = 0.004332755632582322
$ws.Range("F10").Value = 0.08145580589254767
$ws.Range("J10").Value = 0.1377816291161179
$ws.Range("O10").Value = 0.02079722703639515
$ws.Range("Q10").Value = 0.1793760831889082
$ws.Range("R10").Value = 0.07538994800693241
$ws.Range("S10").Value = 0.3500866551126516
$ws.Range("G11").Value = 0.1551020408163265
$ws.Range("J11").Value = 0.09795918367346938
$ws.Range("K11").Value = 0.1918367346938775
$ws.Range("L11").Value = 0.5387755102040817
$ws.Range("S11").Value = 0.0163265306122449
$ws.Range("G12").Value = 0.7913669064748201
$ws.Range("J12").Value = 0.1438848920863309
$ws.Range("K12").Value = 0.007194244604316547
$ws.Range("L12").Value = 0.04316546762589928
$ws.Range("S12").Value = 0.01438848920863309
$ws.Range("G13").Value = 0.6739130434782609
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.06521739130434782
$ws.Range("F15").Value = 0.02040816326530612
$ws.Range("H15").Value = 0.2091836734693878
$ws.Range("I15").Value = 0.08673469387755102
$ws.Range("J15").Value = 0.3418367346938775
$ws.Range("K15").Value = 0.03571428571428571
$ws.Range("M15").Value = 0.00510204081632653
$ws.Range("N15").Value = 0.01020408163265306
$ws.Range("O15").Value = 0.04081632653061224
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.02777777777777778
$ws.Range("H16").Value = 0.1888888888888889
$ws.Range("I16").Value = 0.06111111111111111
$ws.Range("J16").Value = 0.4166666666666667
$ws.Range("K16").Value = 0.1
$ws.Range("M16").Value = 0.02777777777777778
$ws.Range("O16").Value = 0.08888888888888889
$ws.Range("S16").Value = 0.08888888888888889
$ws.Range("F17").Value = 0.01680672268907563
$ws.Range("H17").Value = 0.2240896358543417
$ws.Range("I17").Value = 0.04761904761904762
$ws.Range("J17").Value = 0.4061624649859944
$ws.Range("K17").Value = 0.1372549019607843
$ws.Range("M17").Value = 0.0196078431372549
$ws.Range("O17").Value = 0.06162464985994398
$ws.Range("S17").Value = 0.08683473389355742
$ws.Range("F18").Value = 0.02484472049689441
$ws.Range("H18").Value = 0.1863354037267081
$ws.Range("I18").Value = 0.06211180124223602
$ws.Range("J18").Value = 0.4968944099378882
$ws.Range("K18").Value = 0.09937888198757763
$ws.Range("M18").Value = 0.01863354037267081
$ws.Range("O18").Value = 0.06211180124223602
$ws.Range("S18").Value = 0.04968944099378882
$ws.Range("F19").Value = 0.0215962441314554
$ws.Range("H19").Value = 0.2262910798122066
$ws.Range("I19").Value = 0.05727699530516432
$ws.Range("J19").Value = 0.3849765258215962
$ws.Range("K19").Value = 0.09483568075117371
$ws.Range("M19").Value = 0.03098591549295775
$ws.Range("O19").Value = 0.07042253521126761
$ws.Range("S19").Value = 0.1136150234741784
